# chore: simulator full-month coverage, persist logs, fix employees

$wb = $excel.ActiveWorkbook

$wsTime = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# --- Fix client names (appears on both the timesheet and the schema export) ---
$wsTime.Range("B2").Value = "Markfield"
$wsTime.Range("B3").Value = "Leixner/Smith"

$wsSchema.Range("D2").Value = "Markfield"
$wsSchema.Range("D3").Value = "Leixner/Smith"

# --- Employee ID correction ---
$wsSchema.Range("B2").Value = "emp_qhpjptqm"
$wsSchema.Range("B3").Value = "emp_qhpjptqm"

# --- Weekly Timesheet: rates / totals now populated by simulator ---
$wsTime.Range("E2").Value = 100
$wsTime.Range("F2").Value = 2000

$wsTime.Range("E3").Value = 100
$wsTime.Range("F3").Value = 2000

$wsTime.Range("F5").Value = 4000   # SUBTOTAL
$wsTime.Range("F8").Value = 4000   # HOURLY SUBTOTAL
$wsTime.Range("F10").Value = 4000  # GRAND TOTAL

# --- Jason Schema: mirrored rate / total columns ---
$wsSchema.Range("F2").Value = 100
$wsSchema.Range("G2").Value = 2000

$wsSchema.Range("F3").Value = 100
$wsSchema.Range("G3").Value = 2000
